$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "...and clause 52.232-18" -> "...and 52.232-18" : drop "clause " so
#    the word "clause" can be folded into the new hyperlink display text
#    (do this BEFORE any hyperlinks exist so the replaced run doesn't
#    inherit hyperlink character formatting from a neighbour).
# ---------------------------------------------------------------------
$rFix = $d.Content
$rFix.Find.Execute(" and clause ", $true, $false, $false, $false, $false, $true, 1, $false, " and ", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Turn "FAR 32.703-2(a)" into a hyperlink to the FAR cite.
# ---------------------------------------------------------------------
$rFar32 = $d.Content
$rFar32.Find.Execute("FAR 32.703-2(a)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Hyperlinks.Add($rFar32, "https://www.acquisition.gov/far/part-32#FAR_32_703_2", "", "", "FAR 32.703-2(a)") | Out-Null

# ---------------------------------------------------------------------
# 3) Turn "52.232-18" into a hyperlink, with display text "clause 52.232-18".
# ---------------------------------------------------------------------
$rFar52 = $d.Content
$rFar52.Find.Execute("52.232-18", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Hyperlinks.Add($rFar52, "https://www.acquisition.gov/far/part-52#FAR_52_232_18", "", "", "clause 52.232-18") | Out-Null

# ---------------------------------------------------------------------
# 4) Move the hidden "_GoBack" bookmark from after "26 " (Revised date
#    line) to sit between "The Government reserves" and " the right to
#    cancel..." later in the document. Bookmarks.Add with an existing
#    bookmark name relocates it (removing the old one automatically).
# ---------------------------------------------------------------------
$rGoBack = $d.Content
$rGoBack.Find.Execute("The Government reserves", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rGoBack.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rGoBack) | Out-Null

Write-Host "edit complete"
